# Automatic update of files.
#
# Rows 4 and 5 swap almost all of their field values with each other
# (the two sightings were re-ordered), while column B
# (Taxonsorteringsordning) receives independent new values on both rows.
# Rows 10 and 12 only get their column B value bumped by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (becomes the former row 5 content, plus a new B value) ---
$ws.Range("A4").Value = 131258596
$ws.Range("B4").Value = 91814
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 1202
$ws.Range("F4").Value = "Ullticka"
$ws.Range("G4").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H4").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P4").Value = "Östmossen, Upl"
$ws.Range("Q4").Value = 661129
$ws.Range("R4").Value = 6660729
$ws.Range("Z4").Value = "10:32"
$ws.Range("AB4").Value = "10:32"

# --- Row 5 (becomes the former row 4 content, plus a new B value) ---
$ws.Range("A5").Value = 131258602
$ws.Range("B5").Value = 92536
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 3298
$ws.Range("F5").Value = "Trådticka"
$ws.Range("G5").Value = "Climacocystis borealis"
$ws.Range("H5").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("P5").Value = "Vitmossen, Upl"
$ws.Range("Q5").Value = 661212
$ws.Range("R5").Value = 6660675
$ws.Range("Z5").Value = "09:39"
$ws.Range("AB5").Value = "09:39"

# --- Row 10: only Taxonsorteringsordning changes ---
$ws.Range("B10").Value = 92273

# --- Row 12: only Taxonsorteringsordning changes ---
$ws.Range("B12").Value = 91814
